# Adds BBRB-requested OBI/UO terms (tracker items #823,#824,#825,#834,#835)
# to Sheet1 of the ontoDog input workbook: 29 new rows (150 -> 178 in 1-based
# Excel terms, i.e. sheet rows 149-177), each holding
#   A = source ontology term IRI (hyperlinked for the "Include in View" rows
#       that carry the blue/underlined hyperlink look, same as existing rows)
#   B = source ontology term label
#   C = "y"  (Include in View)
#   E = "yes" (include all children) for the one row that needs it
#
# A handful of B-column cells (149-152) were pasted in with an Arial 10pt
# font (matching other "pasted" rows earlier in the sheet) and one (159)
# has wrap-text turned on - both reproduced below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @{ Row=149; A="http://purl.obolibrary.org/obo/OBI_1000024"; B="cell pellet"; C="y"; E=""; Link=1; ArialB=1; WrapB=0 },
    @{ Row=150; A="http://purl.obolibrary.org/obo/OBI_0001580"; B="material transport service"; C="y"; E=""; Link=1; ArialB=1; WrapB=0 },
    @{ Row=151; A="http://purl.obolibrary.org/obo/OBI_0001173"; B="service"; C="y"; E=""; Link=1; ArialB=1; WrapB=0 },
    @{ Row=152; A="http://purl.obolibrary.org/obo/OBI_0000947"; B="service provider role"; C="y"; E=""; Link=1; ArialB=1; WrapB=0 },
    @{ Row=153; A="http://purl.obolibrary.org/obo/UO_0000195"; B="degree Fahrenheit"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=154; A="http://purl.obolibrary.org/obo/OBI_0002136"; B="RNA Integrity Number calculation"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=155; A="http://purl.obolibrary.org/obo/OBI_0002137"; B="RNA Integrity Number value specification"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=156; A="http://purl.obolibrary.org/obo/OBI_0002145"; B="antigen specific antibodies assay"; C="y"; E="yes"; Link=0; ArialB=0; WrapB=0 },
    @{ Row=157; A="http://purl.obolibrary.org/obo/OBI_0002155"; B="venereal disease research laboratory test"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=158; A="http://purl.obolibrary.org/obo/OBI_0002156"; B="rapid plasma reagin test"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=159; A="http://purl.obolibrary.org/obo/OBI_0002157"; B="HBV surface antigen test"; C="y"; E=""; Link=0; ArialB=0; WrapB=1 },
    @{ Row=160; A="http://purl.obolibrary.org/obo/OBI_0002158"; B="HIV-1 nucleic acid testing"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=161; A="http://purl.obolibrary.org/obo/OBI_0002159"; B="HCV nucleic acid testing"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=162; A="http://purl.obolibrary.org/obo/OBI_0002138"; B="temperature value specification"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=163; A="http://purl.obolibrary.org/obo/OBI_0002139"; B="volume value specification"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=164; A="http://purl.obolibrary.org/obo/OBI_0002140"; B="temperature measurement assay"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=165; A="http://purl.obolibrary.org/obo/OBI_0002141"; B="volume measurement assay"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=166; A="http://purl.obolibrary.org/obo/OBI_0002543"; B="atrial appendage specimen"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=167; A="http://purl.obolibrary.org/obo/OBI_0002544"; B="esophagogastric junction specimen"; C="y"; E=""; Link=0; ArialB=0; WrapB=0 },
    @{ Row=168; A="http://purl.obolibrary.org/obo/OBI_0002545"; B="ileum specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=169; A="http://purl.obolibrary.org/obo/OBI_0002546"; B="liver specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=170; A="http://purl.obolibrary.org/obo/OBI_0002547"; B="minor salivary gland specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=171; A="http://purl.obolibrary.org/obo/OBI_0002548"; B="omentum specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=172; A="http://purl.obolibrary.org/obo/OBI_0002549"; B="ovary specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=173; A="http://purl.obolibrary.org/obo/OBI_0002550"; B="sigmoid colon specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=174; A="http://purl.obolibrary.org/obo/OBI_0002551"; B="suprapubic skin specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=175; A="http://purl.obolibrary.org/obo/OBI_0002552"; B="testis specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=176; A="http://purl.obolibrary.org/obo/OBI_0002553"; B="uterus specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 },
    @{ Row=177; A="http://purl.obolibrary.org/obo/OBI_0002554"; B="vagina specimen"; C="y"; E=""; Link=1; ArialB=0; WrapB=0 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r.A

    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellB.Value = $r.B

    $ws.Cells.Item($rowNum, 3).Value = $r.C

    if ($r.E -ne "") {
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    }

    if ($r.Link -eq 1) {
        $ws.Hyperlinks.Add($cellA, $r.A) | Out-Null
    }

    if ($r.ArialB -eq 1) {
        $cellB.Font.Name = "Arial"
        $cellB.Font.Size = 10
    }

    if ($r.WrapB -eq 1) {
        $cellB.WrapText = $true
    }
}

# Leave the sheet scrolled/selected near the newly-added rows, mirroring the
# author's final cursor position.
$ws.Range("A178").Select()
